# Covariance constraints implemented, but not yet tested for bugs
#
# The diff changes row 3's operator ("eq / ineq") from "=" to "<" and
# updates the active selection on the sheet from B6 to G8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D3 holds the "eq / ineq" operator for the second mean-view constraint.
# Flip it from "=" to "<". Plain Value assignment resets the cell's
# quote-prefix formatting (the column is typed with a leading apostrophe
# so operators like "=" aren't parsed as formulas), so restore the
# original look by copying D2's format (which uses the same style) back
# onto D3 after the value change.
$ws.Range("D3").Value = "<"
$ws.Range("D2").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Move the active selection on the sheet from B6 to G8.
$ws.Range("G8").Select()
